# HeatingElementCostEstimate.xlsx corrections
# - Coil length and coil length (mm vs m) merged into a single row, computed directly in meters
# - "heating plat thickness" -> "heating plate thickness"
# - alum density duplicated (g vs kg) merged into a single row
# - "controler cost" -> "controller cost"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the separate "coil length (m)" row; row 10 will now compute
#        the coil length directly in meters.
$ws.Rows(11).Delete()

# --- 2) Remove the now-redundant "alum density (g/cm^3)" row; the surviving
#        "alum density" row (now row 24) is rewritten below to embed the
#        constant directly.
$ws.Rows(24).Delete()

# --- Row 10: Coil length, now expressed directly in meters ---
$ws.Range("B10").Value = "m"
$ws.Range("C10").Formula = "=(2*3.14159*(89/2))/1000"
$ws.Range("D10").Value = "(2 * pi * (nalgene d)/2)/1000 mm/m"
$ws.Rows(10).RowHeight = 30

# --- Row 16: cross section now references the coil length directly in C10 ---
$ws.Range("C16").Formula = "=((C10*C15))/C14"

# --- Row 17: coil volume now references the coil length directly in C10 ---
$ws.Range("C17").Formula = "=C16*C10"

# --- Row 22: fix spelling "heating plat thickness" -> "heating plate thickness" ---
$ws.Range("A22").Value = "heating plate thickness"

# --- Row 24: alum density, consolidated to a single row (kg/cm^3, embeds 2.7 directly) ---
$ws.Range("A24").Value = "alum density"
$ws.Range("B24").Value = "kg/cm^3"
$ws.Range("C24").Formula = "=2.7*1000"
$ws.Range("D24").Value = "g/cm^3*1000 = kg/m^3"

# --- Row 32: fix spelling "controler cost" -> "controller cost" ---
$ws.Range("A32").Value = "controller cost"

# --- Sheet view: reset top-left cell / selection ---
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
